$d = $word.ActiveDocument

# Locate the "Docente(s) Responsável(eis)" heading paragraph.
$targetIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    if ($para.Range.Text -like "Docente(s)*") {
        $targetIndex = $i
    }
}

if ($targetIndex -eq 0) {
    throw "Could not find 'Docente(s) Responsável(eis)' paragraph"
}

$target = $d.Paragraphs.Item($targetIndex)

# Insert a brand-new empty paragraph right after it.
$target.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs.Item($targetIndex + 1)

# Fill the new paragraph with the bullet-list markup: a "ListBullet"
# styled paragraph containing two runs of text separated by a manual
# line break (<w:br/>), matching the two docentes.
$xmlFragment = '<?xml version="1.0" standalone="yes"?><?mso-application progid="Word.Document"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="ListBullet"/></w:pPr><w:r><w:t>3577649 - Carlos Angelo Nunes</w:t><w:br/></w:r><w:r><w:t>3586455 - Cassius Olivio Figueiredo Terra Ruchert</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$newPara.Range.InsertXML($xmlFragment)
